# Update projection outputs and fix termination compensation proration logic
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 102;  "C2" = 102;  "D2" = 102
    "F2" = 1
    "G2" = 0.09127332586241999
    "H2" = 0.09127332586241999
    "I2" = 544219.5395527922
    "J2" = 209495.769977396
    "L2" = 209495.769977396
    "M2" = 753715.3095301883
    "N2" = 10316742.6388
    "O2" = 9909001.708699998
    "P2" = 0.02030638713323217
    "Q2" = 0.02114196526916137

    "B3" = 106;  "C3" = 106;  "D3" = 106
    "F3" = 1
    "G3" = 0.09717788792471758
    "H3" = 0.09717788792471758
    "I3" = 654630.7013274725
    "J3" = 262102.9308707663
    "L3" = 262102.9308707663
    "M3" = 916733.6321982386
    "N3" = 10760730.661864
    "O3" = 10353357.503861
    "P3" = 0.02435735444988493
    "Q3" = 0.02531574233508523

    "B4" = 108;  "C4" = 108;  "D4" = 108
    "G4" = 0.09963069641345117
    "H4" = 0.09963069641345117
    "I4" = 734817.7333005213
    "J4" = 295110.0742635016
    "L4" = 295110.0742635016
    "M4" = 1029927.807564023
    "N4" = 11223740.87971992
    "O4" = 10815396.52697683
    "P4" = 0.02629337913500242
    "Q4" = 0.02728610768245148

    "B5" = 109;  "C5" = 109;  "D5" = 109
    "G5" = 0.1015414746775767
    "H5" = 0.1015414746775767
    "I5" = 781388.0653614923
    "J5" = 315821.2765223843
    "L5" = 315821.2765223843
    "M5" = 1097209.341883876
    "N5" = 11501137.06941152
    "O5" = 11090442.38608613
    "P5" = 0.02746000457314296
    "Q5" = 0.02847688717256292

    "B6" = 110;  "C6" = 110;  "D6" = 110
    "G6" = 0.1022024053100728
    "H6" = 0.1022024053100728
    "I6" = 827577.8738163244
    "J6" = 335854.9980650494
    "L6" = 335854.9980650494
    "M6" = 1163432.871881374
    "N6" = 11926849.94239386
    "O6" = 11512384.41856872
    "P6" = 0.02815957270253366
    "Q6" = 0.02917336546921916
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
